$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark attendance (week 14 / column P) for the rows that attended.
$rows = @(3, 4, 6, 8, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 22)
foreach ($r in $rows) {
    $ws.Range("P$r").Value = 1
}

# Update the active selection to match where the editor last clicked.
$ws.Range("P22").Select()
